# Uganda central resources workbook — update the "Notes" sheet metadata:
# refine the description/source text, add a source link, and split the
# licensing note into a statement + a link to the license text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Insert a new row right after "Source: ..." (row 4) to hold the new
# "Source-link" line; this pushes the rest of the notes down by one row.
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Source-link: http://www.budget.go.ug/"

# Refine the description and source lines.
$ws.Range("A2").Value = "Description: Centrally Raised Revenues"
$ws.Range("A4").Value = "Source: Local Government Budgets - Ministry of Finance, Planning and Economic Development"

# Update the licensing sentence (previously "It is provided on an as-is
# basis under an open-use license."), then insert a new row below it with a
# link to the license text.
$ws.Range("A16").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
